$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.616.31"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "3.497.63"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.12"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.02"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E8").Value = "  +4.65%  "
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.431"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "4.105.03"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.14"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").Value = "66.662.83"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "3.494.14"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.31"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.05"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "393.64"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.94"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.11"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.533"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.18"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.35"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.45"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.05"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.79"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +5.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.25"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.897"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.81"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.65"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.54"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.798.52"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.79"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.89"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.95"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.96"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.852"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("E51").Value = "  +1.88%  "
